$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 282
$ws1.Range("F8").Value  = 7631
$ws1.Range("F11").Value = 16
$ws1.Range("F14").Value = 658
$ws1.Range("F16").Value = 1044
$ws1.Range("F19").Value = 1515
$ws1.Range("F20").Value = 327
$ws1.Range("F21").Value = 6037
$ws1.Range("F22").Value = 33
$ws1.Range("F27").Value = 4185
$ws1.Range("F28").Value = 3836
$ws1.Range("F29").Value = 283
$ws1.Range("F34").Value = 1019
$ws1.Range("F37").Value = 73
$ws1.Range("F40").Value = 189
$ws1.Range("F43").Value = 374
$ws1.Range("F45").Value = 1080
$ws1.Range("F47").Value = 2875
$ws1.Range("F49").Value = 331

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 373
$ws2.Range("F10").Value = 600
$ws2.Range("F20").Value = 148
$ws2.Range("F28").Value = 5124
$ws2.Range("F29").Value = 5124

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value  = 3024
$ws3.Range("F13").Value = 2044
$ws3.Range("F14").Value = 8740

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 282
$ws4.Range("F7").Value  = 3024
$ws4.Range("F8").Value  = 7631
$ws4.Range("F17").Value = 373
$ws4.Range("F18").Value = 600
$ws4.Range("F19").Value = 600
$ws4.Range("F20").Value = 658
$ws4.Range("F22").Value = 1044
$ws4.Range("F27").Value = 1515
$ws4.Range("F28").Value = 327
$ws4.Range("F29").Value = 6037
$ws4.Range("F31").Value = 4185
$ws4.Range("F32").Value = 3836
$ws4.Range("F36").Value = 1019
$ws4.Range("F38").Value = 73
$ws4.Range("F40").Value = 189
$ws4.Range("F42").Value = 374
$ws4.Range("F46").Value = 2876
$ws4.Range("F48").Value = 5124
